# Boost Spirit karma_performance.xlsx — "added one more Karma benchmark"
#
# Adds a new "Sheet1" worksheet (after "Single int") containing a benchmark
# table for sequences of different length, plus a clustered-column chart
# ("Karma sequences") comparing six compilers across sequence lengths 2-9.

$wb = $excel.ActiveWorkbook

# --- Leave a trace of the previous selection on "Single int" (it loses the
#     active-tab/tabSelected flag once the new sheet becomes active) --------
$wsSingleInt = $wb.Worksheets.Item("Single int")
$wsSingleInt.Range("E3").Select()

# --- Add the new worksheet at the end of the tab strip ---------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Sheet1"

# --- Caption ------------------------------------------------------------
$ws.Range("E1").Value = "Benchmarking sequence of different length:"

# --- Column headers (row 4): compiler / tool names --------------------
$ws.Range("F4").Value = "VC8SP1"
$ws.Range("G4").Value = "gcc 4.4.0 (32)"
$ws.Range("H4").Value = "VC++ 10 (32)"
$ws.Range("I4").Value = "Intel 11.1 (32)"
$ws.Range("J4").Value = "gcc 4.4.0 (64)"
$ws.Range("K4").Value = "VC++ 10 (64)"
$ws.Range("L4").Value = "Intel 11.1 (64)"

# --- Data rows (5-12): sequence length (E) + measured times (F-L) ------
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.256
$ws.Range("G5").Value = 0.627
$ws.Range("H5").Value = 0.424
$ws.Range("I5").Value = 0.569
$ws.Range("J5").Value = 0.819
$ws.Range("K5").Value = 0.311
$ws.Range("L5").Value = 0.426

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 0.7
$ws.Range("G6").Value = 0.956
$ws.Range("H6").Value = 0.492
$ws.Range("I6").Value = 0.98
$ws.Range("J6").Value = 1.277
$ws.Range("K6").Value = 0.441
$ws.Range("L6").Value = 0.695

$ws.Range("E7").Value = 4
$ws.Range("F7").Value = 1.044
$ws.Range("G7").Value = 1.298
$ws.Range("H7").Value = 1.217
$ws.Range("I7").Value = 1.448
$ws.Range("J7").Value = 1.472
$ws.Range("K7").Value = 0.768
$ws.Range("L7").Value = 1.019

$ws.Range("E8").Value = 5
$ws.Range("F8").Value = 1.443
$ws.Range("G8").Value = 1.65
$ws.Range("H8").Value = 1.789
$ws.Range("I8").Value = 1.989
$ws.Range("J8").Value = 2.229
$ws.Range("K8").Value = 0.975
$ws.Range("L8").Value = 1.265

$ws.Range("E9").Value = 6
$ws.Range("F9").Value = 1.76
$ws.Range("G9").Value = 1.906
$ws.Range("H9").Value = 1.986
$ws.Range("I9").Value = 2.596
$ws.Range("J9").Value = 2.709
$ws.Range("K9").Value = 1.216
$ws.Range("L9").Value = 1.634

$ws.Range("E10").Value = 7
$ws.Range("F10").Value = 2.214
$ws.Range("G10").Value = 2.478
$ws.Range("H10").Value = 2.513
$ws.Range("I10").Value = 3.242
$ws.Range("J10").Value = 3.205
$ws.Range("K10").Value = 1.634
$ws.Range("L10").Value = 2.008

$ws.Range("E11").Value = 8
$ws.Range("F11").Value = 2.756
$ws.Range("G11").Value = 2.676
$ws.Range("H11").Value = 2.829
$ws.Range("I11").Value = 3.559
$ws.Range("J11").Value = 3.377
$ws.Range("K11").Value = 1.853
$ws.Range("L11").Value = 2.324

$ws.Range("E12").Value = 9
$ws.Range("F12").Value = 3.29
$ws.Range("G12").Value = 2.945
$ws.Range("H12").Value = 3.732
$ws.Range("I12").Value = 4.246
$ws.Range("J12").Value = 3.592
$ws.Range("K12").Value = 2.11
$ws.Range("L12").Value = 2.662

# --- Column widths, roughly matching the sibling benchmark sheets ------
$ws.Columns.Item(6).ColumnWidth = 11.29
$ws.Columns.Item(7).ColumnWidth = 12.14
$ws.Columns.Item(8).ColumnWidth = 11.71
$ws.Columns.Item(9).ColumnWidth = 13.14
$ws.Columns.Item(10).ColumnWidth = 12.14
$ws.Columns.Item(11).ColumnWidth = 11.71
$ws.Columns.Item(12).ColumnWidth = 13.14

# --- Chart: "Karma sequences" clustered column chart --------------------
$chartObj = $ws.ChartObjects().Add(200, 20, 430, 330)
$chart = $chartObj.Chart
$chart.ChartType = 51          # xlColumnClustered
$chart.ChartStyle = 18

$chart.HasTitle = $true
$chart.ChartTitle.Text = "Karma sequences"

$catRange = $ws.Range("E5:E12")

$s1 = $chart.SeriesCollection().NewSeries()
$s1.Name = "=Sheet1!`$G`$4"
$s1.XValues = $catRange
$s1.Values = $ws.Range("G5:G12")

$s2 = $chart.SeriesCollection().NewSeries()
$s2.Name = "=Sheet1!`$H`$4"
$s2.XValues = $catRange
$s2.Values = $ws.Range("H5:H12")

$s3 = $chart.SeriesCollection().NewSeries()
$s3.Name = "=Sheet1!`$I`$4"
$s3.XValues = $catRange
$s3.Values = $ws.Range("I5:I12")

$s4 = $chart.SeriesCollection().NewSeries()
$s4.Name = "=Sheet1!`$J`$4"
$s4.XValues = $catRange
$s4.Values = $ws.Range("J5:J12")

$s5 = $chart.SeriesCollection().NewSeries()
$s5.Name = "=Sheet1!`$K`$4"
$s5.XValues = $catRange
$s5.Values = $ws.Range("K5:K12")

$s6 = $chart.SeriesCollection().NewSeries()
$s6.Name = "=Sheet1!`$L`$4"
$s6.XValues = $catRange
$s6.Values = $ws.Range("L5:L12")

$chart.ChartGroups(1).GapWidth = 152

$catAxis = $chart.Axes(1)
$catAxis.HasTitle = $true
$catAxis.AxisTitle.Text = "Sequence length (elements)"

$valAxis = $chart.Axes(2)
$valAxis.HasTitle = $true
$valAxis.AxisTitle.Text = "Measured time [s]"
$valAxis.HasMajorGridlines = $true

$chart.HasLegend = $true
$chart.Legend.Position = -4107   # xlLegendPositionBottom

# --- Re-select the first data cell on the new sheet, like the author left it
$ws.Range("G5").Select()

Write-Host "Sheet1 added with benchmark data and chart."
